$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.043892325127094
$ws.Cells.Item(2, 4).Value = 1.051176082905937
$ws.Cells.Item(2, 5).Value = 1.051260085145705
$ws.Cells.Item(2, 6).Value = 1.060754753668595
$ws.Cells.Item(2, 9).Value = 1.027899155600938
$ws.Cells.Item(2, 10).Value = 1.048960335985511
$ws.Cells.Item(2, 11).Value = 1.05392831161863
$ws.Cells.Item(2, 12).Value = 1.054012080863353
$ws.Cells.Item(2, 13).Value = 1.06348067176402
$ws.Cells.Item(2, 14).Value = 1.050449979980043

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.045152066864657
$ws.Cells.Item(3, 4).Value = 1.05230253714161
$ws.Cells.Item(3, 5).Value = 1.052369830490029
$ws.Cells.Item(3, 6).Value = 1.061925473886423
$ws.Cells.Item(3, 9).Value = 1.027874080474153
$ws.Cells.Item(3, 10).Value = 1.049865573522689
$ws.Cells.Item(3, 11).Value = 1.054866552268061
$ws.Cells.Item(3, 12).Value = 1.054933672550396
$ws.Cells.Item(3, 13).Value = 1.064464980609586
$ws.Cells.Item(3, 14).Value = 1.051356503058356

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.045967373687923
$ws.Cells.Item(4, 4).Value = 1.053031728201908
$ws.Cells.Item(4, 5).Value = 1.053088207616525
$ws.Cells.Item(4, 6).Value = 1.062682989192161
$ws.Cells.Item(4, 9).Value = 1.027855467289442
$ws.Cells.Item(4, 10).Value = 1.050450997333966
$ws.Cells.Item(4, 11).Value = 1.055473370943869
$ws.Cells.Item(4, 12).Value = 1.055529712462965
$ws.Cells.Item(4, 13).Value = 1.065101296787718
$ws.Cells.Item(4, 14).Value = 1.051942758238594

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.046310171130769
$ws.Cells.Item(5, 4).Value = 1.053338352761538
$ws.Cells.Item(5, 5).Value = 1.053390285410622
$ws.Cells.Item(5, 6).Value = 1.063001444767751
$ws.Cells.Item(5, 9).Value = 1.027847069720865
$ws.Cells.Item(5, 10).Value = 1.050697032464499
$ws.Cells.Item(5, 11).Value = 1.055728409868067
$ws.Cells.Item(5, 12).Value = 1.055780218563471
$ws.Cells.Item(5, 13).Value = 1.065368661604129
$ws.Cells.Item(5, 14).Value = 1.052189142767233

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.046367730872623
$ws.Cells.Item(6, 4).Value = 1.05338984065862
$ws.Cells.Item(6, 5).Value = 1.053441009849267
$ws.Cells.Item(6, 6).Value = 1.06305491465355
$ws.Cells.Item(6, 9).Value = 1.027845626128536
$ws.Cells.Item(6, 10).Value = 1.050738338352934
$ws.Cells.Item(6, 11).Value = 1.055771228080995
$ws.Cells.Item(6, 12).Value = 1.05582227560995
$ws.Cells.Item(6, 13).Value = 1.065413544991993
$ws.Cells.Item(6, 14).Value = 1.052230507314768

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.045971953995039
$ws.Cells.Item(7, 4).Value = 1.05303582504782
$ws.Cells.Item(7, 5).Value = 1.053092243710731
$ws.Cells.Item(7, 6).Value = 1.062687244425503
$ws.Cells.Item(7, 9).Value = 1.027855357331768
$ws.Cells.Item(7, 10).Value = 1.05045428516908
$ws.Cells.Item(7, 11).Value = 1.055476779051335
$ws.Cells.Item(7, 12).Value = 1.055533060008307
$ws.Cells.Item(7, 13).Value = 1.065104869888811
$ws.Cells.Item(7, 14).Value = 1.051946050742811

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.044318026714035
$ws.Cells.Item(8, 4).Value = 1.051556711818948
$ws.Cells.Item(8, 5).Value = 1.051635067616774
$ws.Cells.Item(8, 6).Value = 1.061150408384004
$ws.Cells.Item(8, 9).Value = 1.027891175227185
$ws.Cells.Item(8, 10).Value = 1.049266333147558
$ws.Cells.Item(8, 11).Value = 1.054245454069522
$ws.Cells.Item(8, 12).Value = 1.054323597884408
$ws.Cells.Item(8, 13).Value = 1.063813447576045
$ws.Cells.Item(8, 14).Value = 1.050756411693158

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.041404812874007
$ws.Cells.Item(9, 4).Value = 1.048952582998858
$ws.Cells.Item(9, 5).Value = 1.049069581362021
$ws.Cells.Item(9, 6).Value = 1.05844213052052
$ws.Cells.Item(9, 9).Value = 1.027936042409227
$ws.Cells.Item(9, 10).Value = 1.047170478584806
$ws.Cells.Item(9, 11).Value = 1.052073478754376
$ws.Cells.Item(9, 12).Value = 1.052190105743107
$ws.Cells.Item(9, 13).Value = 1.061533196148025
$ws.Cells.Item(9, 14).Value = 1.048657580776529

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.039463371585489
$ws.Cells.Item(10, 4).Value = 1.04721795789742
$ws.Cells.Item(10, 5).Value = 1.047360713227455
$ws.Cells.Item(10, 6).Value = 1.056636446214352
$ws.Cells.Item(10, 9).Value = 1.027953737677185
$ws.Cells.Item(10, 10).Value = 1.04577148322352
$ws.Cells.Item(10, 11).Value = 1.050623949075276
$ws.Cells.Item(10, 12).Value = 1.050766206399301
$ws.Cells.Item(10, 13).Value = 1.060009898031397
$ws.Cells.Item(10, 14).Value = 1.047256598681363

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.038622845617691
$ws.Cells.Item(11, 4).Value = 1.046467178929251
$ws.Cells.Item(11, 5).Value = 1.04662108713725
$ws.Cells.Item(11, 6).Value = 1.055854514085219
$ws.Cells.Item(11, 9).Value = 1.02795851336104
$ws.Cells.Item(11, 10).Value = 1.045165272122337
$ws.Cells.Item(11, 11).Value = 1.049995907762216
$ws.Cells.Item(11, 12).Value = 1.050149258038087
$ws.Cells.Item(11, 13).Value = 1.059349541231887
$ws.Cells.Item(11, 14).Value = 1.046649526690884

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.038310654118596
$ws.Cells.Item(12, 4).Value = 1.046188353764212
$ws.Cells.Item(12, 5).Value = 1.046346404720626
$ws.Cells.Item(12, 6).Value = 1.055564059703944
$ws.Cells.Item(12, 9).Value = 1.027959854441096
$ws.Cells.Item(12, 10).Value = 1.04494003143894
$ws.Cells.Item(12, 11).Value = 1.049762566199346
$ws.Cells.Item(12, 12).Value = 1.049920036177089
$ws.Cells.Item(12, 13).Value = 1.059104140401754
$ws.Cells.Item(12, 14).Value = 1.04642396613988

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.038377619511988
$ws.Cells.Item(13, 4).Value = 1.046248160576435
$ws.Cells.Item(13, 5).Value = 1.046405322894638
$ws.Cells.Item(13, 6).Value = 1.055626363624096
$ws.Cells.Item(13, 9).Value = 1.027959586347597
$ws.Cells.Item(13, 10).Value = 1.04498834937715
$ws.Cells.Item(13, 11).Value = 1.049812621447654
$ws.Cells.Item(13, 12).Value = 1.049969207766728
$ws.Cells.Item(13, 13).Value = 1.059156784927913
$ws.Cells.Item(13, 14).Value = 1.046472352695105

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.038597039432583
$ws.Cells.Item(14, 4).Value = 1.046444130191081
$ws.Cells.Item(14, 5).Value = 1.046598380837996
$ws.Cells.Item(14, 6).Value = 1.055830505245422
$ws.Cells.Item(14, 9).Value = 1.02795863303975
$ws.Cells.Item(14, 10).Value = 1.04514665502896
$ws.Cells.Item(14, 11).Value = 1.049976620887251
$ws.Cells.Item(14, 12).Value = 1.050130311713191
$ws.Cells.Item(14, 13).Value = 1.059329258672312
$ws.Cells.Item(14, 14).Value = 1.046630883159098

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.038732233503999
$ws.Cells.Item(15, 4).Value = 1.04656487980394
$ws.Cells.Item(15, 5).Value = 1.04671733648619
$ws.Cells.Item(15, 6).Value = 1.055956282280616
$ws.Cells.Item(15, 9).Value = 1.027957988347182
$ws.Cells.Item(15, 10).Value = 1.045244183464805
$ws.Cells.Item(15, 11).Value = 1.050077658494588
$ws.Cells.Item(15, 12).Value = 1.050229565218772
$ws.Cells.Item(15, 13).Value = 1.059435510183572
$ws.Cells.Item(15, 14).Value = 1.046728550096508

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.039519157079385
$ws.Cells.Item(16, 4).Value = 1.047267791366468
$ws.Cells.Item(16, 5).Value = 1.047409806522469
$ws.Cells.Item(16, 6).Value = 1.056688339137739
$ws.Cells.Item(16, 9).Value = 1.027953359987855
$ws.Cells.Item(16, 10).Value = 1.04581170614644
$ws.Cells.Item(16, 11).Value = 1.050665621890557
$ws.Cells.Item(16, 12).Value = 1.050807142908445
$ws.Cells.Item(16, 13).Value = 1.060053707673818
$ws.Cells.Item(16, 14).Value = 1.047296878725448

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.04001280680573
$ws.Cells.Item(17, 4).Value = 1.047708795369358
$ws.Cells.Item(17, 5).Value = 1.047844260873858
$ws.Cells.Item(17, 6).Value = 1.057147522501093
$ws.Cells.Item(17, 9).Value = 1.027949684543335
$ws.Cells.Item(17, 10).Value = 1.046167580282625
$ws.Cells.Item(17, 11).Value = 1.051034331799684
$ws.Cells.Item(17, 12).Value = 1.051169336801134
$ws.Cells.Item(17, 13).Value = 1.060441282855114
$ws.Cells.Item(17, 14).Value = 1.047653258243725

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.040300757134736
$ws.Cells.Item(18, 4).Value = 1.047966057002947
$ws.Cells.Item(18, 5).Value = 1.048097702195357
$ws.Cells.Item(18, 6).Value = 1.057415350782008
$ws.Cells.Item(18, 9).Value = 1.027947262159446
$ws.Cells.Item(18, 10).Value = 1.046375113578911
$ws.Cells.Item(18, 11).Value = 1.051249356945793
$ws.Cells.Item(18, 12).Value = 1.051380560763919
$ws.Cells.Item(18, 13).Value = 1.060667275765827
$ws.Cells.Item(18, 14).Value = 1.047861086261097

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.040398942998016
$ws.Cells.Item(19, 4).Value = 1.048053781945963
$ws.Cells.Item(19, 5).Value = 1.048184124506306
$ws.Cells.Item(19, 6).Value = 1.057506672386863
$ws.Cells.Item(19, 9).Value = 1.027946388914831
$ws.Cells.Item(19, 10).Value = 1.046445869972014
$ws.Cells.Item(19, 11).Value = 1.051322668704465
$ws.Cells.Item(19, 12).Value = 1.051452576330494
$ws.Cells.Item(19, 13).Value = 1.060744321130741
$ws.Cells.Item(19, 14).Value = 1.047931943136394

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.039959841557024
$ws.Cells.Item(20, 4).Value = 1.047661476569679
$ws.Cells.Item(20, 5).Value = 1.047797644790344
$ws.Cells.Item(20, 6).Value = 1.057098257060946
$ws.Cells.Item(20, 9).Value = 1.027950107690478
$ws.Cells.Item(20, 10).Value = 1.04612940271685
$ws.Cells.Item(20, 11).Value = 1.050994776563152
$ws.Cells.Item(20, 12).Value = 1.051130480722889
$ws.Cells.Item(20, 13).Value = 1.06039970727857
$ws.Cells.Item(20, 14).Value = 1.047615026461428

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.038532425277425
$ws.Cells.Item(21, 4).Value = 1.046386420721044
$ws.Cells.Item(21, 5).Value = 1.046541528779912
$ws.Cells.Item(21, 6).Value = 1.055770390910393
$ws.Cells.Item(21, 9).Value = 1.027958925706458
$ws.Cells.Item(21, 10).Value = 1.045100039816711
$ws.Cells.Item(21, 11).Value = 1.049928328779388
$ws.Cells.Item(21, 12).Value = 1.050082872275128
$ws.Cells.Item(21, 13).Value = 1.059278472633046
$ws.Cells.Item(21, 14).Value = 1.0465842017479

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.03763505073521
$ws.Cells.Item(22, 4).Value = 1.045585017103554
$ws.Cells.Item(22, 5).Value = 1.045752033852957
$ws.Cells.Item(22, 6).Value = 1.054935450496796
$ws.Cells.Item(22, 9).Value = 1.027961966178107
$ws.Cells.Item(22, 10).Value = 1.044452450884162
$ws.Cells.Item(22, 11).Value = 1.049257468477535
$ws.Cells.Item(22, 12).Value = 1.049423852908053
$ws.Cells.Item(22, 13).Value = 1.058572842595227
$ws.Cells.Item(22, 14).Value = 1.045935693164788

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.038110757403342
$ws.Cells.Item(23, 4).Value = 1.046009830723481
$ws.Cells.Item(23, 5).Value = 1.046170534379576
$ws.Cells.Item(23, 6).Value = 1.05537807410726
$ws.Cells.Item(23, 9).Value = 1.027960591417237
$ws.Cells.Item(23, 10).Value = 1.044795787191923
$ws.Cells.Item(23, 11).Value = 1.04961313704523
$ws.Cells.Item(23, 12).Value = 1.049773244741204
$ws.Cells.Item(23, 13).Value = 1.058946973806699
$ws.Cells.Item(23, 14).Value = 1.046279517049488

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.039983774233538
$ws.Cells.Item(24, 4).Value = 1.047682857801351
$ws.Cells.Item(24, 5).Value = 1.04781870849487
$ws.Cells.Item(24, 6).Value = 1.057120518010208
$ws.Cells.Item(24, 9).Value = 1.027949917349334
$ws.Cells.Item(24, 10).Value = 1.046146653645544
$ws.Cells.Item(24, 11).Value = 1.05101264998665
$ws.Cells.Item(24, 12).Value = 1.051148038228404
$ws.Cells.Item(24, 13).Value = 1.060418493718115
$ws.Cells.Item(24, 14).Value = 1.047632301888419

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.042157816278689
$ws.Cells.Item(25, 4).Value = 1.049625551039761
$ws.Cells.Item(25, 5).Value = 1.049732560525132
$ws.Cells.Item(25, 6).Value = 1.05914231099548
$ws.Cells.Item(25, 9).Value = 1.027926600194826
$ws.Cells.Item(25, 10).Value = 1.047712613062052
$ws.Cells.Item(25, 11).Value = 1.052635255177078
$ws.Cells.Item(25, 12).Value = 1.052741937996899
$ws.Cells.Item(25, 13).Value = 1.062123244795018
$ws.Cells.Item(25, 14).Value = 1.049200485146917

$wb.Save()